# Modify the FA (fatty acid) table on the "all" worksheet:
#   - add a new column AQ ("g extracted") with per-sample weights
#   - rework the "total FAs ug/g" (AP) formulas to divide by both the
#     %recovery-normalised total (AO) and the newly added "g extracted" (AQ)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all")

# --- header for the new column AQ -----------------------------------------
$ws.Range("AQ3").Value = "g extracted"

# --- new "g extracted" values, one per sample (rows 4-9) -------------------
$ws.Range("AQ4").Value = 3.53
$ws.Range("AQ5").Value = 3.52
$ws.Range("AQ6").Value = 3.49
$ws.Range("AQ7").Value = 3.51
$ws.Range("AQ8").Value = 3.57
$ws.Range("AQ9").Value = 3.52

# --- updated "total FAs ug/g" formulas --------------------------------------
# Row 4 keeps its own (non-shared) formula.
$ws.Range("AP4").Formula = "= SUM(W4:AL4)/AO4/AQ4/1000"
$ws.Range("AP4").Style = "Normal"

# Rows 5-9 share one formula (relative references adjust per row).
$ws.Range("AP5:AP9").Formula = "= SUM(W5:AL5)/AO5/AQ5/1000"
$ws.Range("AP5:AP9").Style = "Normal"

# Row 10 (blank control) has no %recovery or "g extracted" value, so the
# total is simply divided by 1000.
$ws.Range("AP10").Formula = "= SUM(W10:AL10)/1000"
$ws.Range("AP10").Style = "Normal"

# --- keep the active selection where the author left it --------------------
$ws.Range("AP9").Select()
